$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 351.57144
$ws.Range("I6").Value = 300.25
$ws.Range("J6").Value = 420
$ws.Range("K6").Value = 900.75
$ws.Range("L6").Value = 1260
$ws.Range("M6").Value = -788.75
$ws.Range("N6").Value = -1484

$ws.Range("H51").Value = 44443.145
$ws.Range("J51").Value = 22220
$ws.Range("L51").Value = 22220
$ws.Range("N51").Value = -23188

$ws.Range("H76").Value = 11119017
$ws.Range("I76").Value = 10122.714
$ws.Range("J76").Value = 20839300
$ws.Range("K76").Value = 10122.714
$ws.Range("L76").Value = 20839300
$ws.Range("M76").Value = -9807.714
$ws.Range("N76").Value = -20839930

$ws.Range("H79").Value = 11119017
$ws.Range("I79").Value = 10122.714
$ws.Range("J79").Value = 20839300
$ws.Range("K79").Value = 10122.714
$ws.Range("L79").Value = 20839300
$ws.Range("M79").Value = -9030.714
$ws.Range("N79").Value = -20841484

$ws.Range("H107").Value = 640
$ws.Range("I107").Value = 480
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 480
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1440
$ws.Range("N107").Value = -4640

$ws.Range("H111").Value = 1088.6666
$ws.Range("I111").Value = 733.3333
$ws.Range("J111").Value = 1444
$ws.Range("K111").Value = 2199.9999
$ws.Range("L111").Value = 4332
$ws.Range("M111").Value = 867.0001000000002
$ws.Range("N111").Value = -10466

$ws.Range("H139").Value = 42900.5
$ws.Range("J139").Value = 42900.5
$ws.Range("L139").Value = 42900.5
$ws.Range("N139").Value = -53180.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12989418
$ws.Range("I32").Value = 1878.197
$ws.Range("K32").Value = 1878.197
$ws.Range("M32").Value = -1591.197

$ws.Range("H74").Value = 2132.261
$ws.Range("I74").Value = 2315.2666
$ws.Range("J74").Value = 1789.125
$ws.Range("K74").Value = 2315.2666
$ws.Range("L74").Value = 1789.125
$ws.Range("M74").Value = -1441.2666
$ws.Range("N74").Value = -3537.125

$ws.Range("H77").Value = 2132.261
$ws.Range("I77").Value = 2315.2666
$ws.Range("J77").Value = 1789.125
$ws.Range("K77").Value = 11576.333
$ws.Range("L77").Value = 8945.625
$ws.Range("M77").Value = -7208.332999999999
$ws.Range("N77").Value = -17681.625

$ws.Range("H101").Value = 36301
$ws.Range("J101").Value = 36301
$ws.Range("L101").Value = 36301
$ws.Range("N101").Value = -42791

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 23000
$ws.Range("J35").Value = 23000
$ws.Range("L35").Value = 23000
$ws.Range("N35").Value = -23620

$ws.Range("H102").Value = 10556
$ws.Range("I102").Value = 10556
$ws.Range("K102").Value = 10556
$ws.Range("M102").Value = -7311

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 434.44116
$ws.Range("I107").Value = 346.86957
$ws.Range("J107").Value = 617.5454999999999
$ws.Range("K107").Value = 346.86957
$ws.Range("L107").Value = 617.5454999999999
$ws.Range("M107").Value = 1573.13043
$ws.Range("N107").Value = -4457.5455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 36462110
$ws.Range("I5").Value = 55555924
$ws.Range("K5").Value = 166667772
$ws.Range("M5").Value = -166667660

$ws.Range("H113").Value = 8238119.5
$ws.Range("I113").Value = 5952965
$ws.Range("J113").Value = 10370930
$ws.Range("K113").Value = 17858895
$ws.Range("L113").Value = 31112790
$ws.Range("M113").Value = -17856725
$ws.Range("N113").Value = -31117130

$ws.Range("H132").Value = 5911.4546
$ws.Range("J132").Value = 6443.2
$ws.Range("L132").Value = 57988.8
$ws.Range("N132").Value = -63048.8

$ws.Range("H135").Value = 36462110
$ws.Range("I135").Value = 55555924
$ws.Range("K135").Value = 500003316
$ws.Range("M135").Value = -500000781

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12504605
$ws.Range("I80").Value = 6975
$ws.Range("J80").Value = 25002236
$ws.Range("K80").Value = 6975
$ws.Range("L80").Value = 25002236
$ws.Range("M80").Value = -5977
$ws.Range("N80").Value = -25004232

$ws.Range("H83").Value = 12504605
$ws.Range("I83").Value = 6975
$ws.Range("J83").Value = 25002236
$ws.Range("K83").Value = 34875
$ws.Range("L83").Value = 125011180
$ws.Range("M83").Value = -29883
$ws.Range("N83").Value = -125021164

$ws.Range("H132").Value = 9320.134
$ws.Range("I132").Value = 1687.5
$ws.Range("J132").Value = 12095.637
$ws.Range("K132").Value = 5062.5
$ws.Range("L132").Value = 36286.911
$ws.Range("M132").Value = -2532.5
$ws.Range("N132").Value = -41346.911

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2659
$ws.Range("I22").Value = 3666.6667
$ws.Range("J22").Value = 2227.1428
$ws.Range("K22").Value = 3666.6667
$ws.Range("L22").Value = 2227.1428
$ws.Range("M22").Value = -3371.6667
$ws.Range("N22").Value = -2817.1428

$ws.Range("H27").Value = 2659
$ws.Range("I27").Value = 3666.6667
$ws.Range("J27").Value = 2227.1428
$ws.Range("K27").Value = 3666.6667
$ws.Range("L27").Value = 2227.1428
$ws.Range("M27").Value = -3559.6667
$ws.Range("N27").Value = -2441.1428

$ws.Range("H132").Value = 42338516
$ws.Range("I132").Value = 71430840
$ws.Range("J132").Value = 22407.908
$ws.Range("K132").Value = 214292520
$ws.Range("L132").Value = 67223.724
$ws.Range("M132").Value = -214289990
$ws.Range("N132").Value = -72283.724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1275.15
$ws.Range("I107").Value = 1478.5454
$ws.Range("J107").Value = 1026.5555
$ws.Range("K107").Value = 4435.6362
$ws.Range("L107").Value = 3079.6665
$ws.Range("M107").Value = -2515.6362
$ws.Range("N107").Value = -6919.666499999999

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws.Range("H139").Value = 71050
$ws.Range("J139").Value = 71050
$ws.Range("L139").Value = 71050
$ws.Range("N139").Value = -81330

